$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "COMP-THC-6106"
$ws.Range("J3").Value = "COMP-THC-6106"
$ws.Range("J4").Value = "COMP-THC-6106"
$ws.Range("J5").Value = "THC-6106"
$ws.Range("J6").Value = "THC-6107"
$ws.Range("J7").Value = "THC-6108"
$ws.Range("J8").Value = "COMP-THC-6109"
$ws.Range("J9").Value = "COMP-THC-6109"
$ws.Range("J10").Value = "COMP-THC-6109"
$ws.Range("J11").Value = "THC-6109"
$ws.Range("J12").Value = "THC-6110"
$ws.Range("J13").Value = "THC-6111"
$ws.Range("J14").Value = "COMP-THC-6112"
$ws.Range("J15").Value = "COMP-THC-6112"
$ws.Range("J16").Value = "COMP-THC-6112"
$ws.Range("J17").Value = "THC-6112"
$ws.Range("J18").Value = "THC-6113"
$ws.Range("J19").Value = "THC-6114"
$ws.Range("J20").Value = "THC-6114"
$ws.Range("J21").Value = "THC-6114"
$ws.Range("J22").Value = "COMP-THC-6114"
$ws.Range("J23").Value = "COMP-THC-6115"
$ws.Range("J24").Value = "THC-6116"
$ws.Range("J25").Value = "COMP-THC-6117"
$ws.Range("J26").Value = "THC-6117"
$ws.Range("J27").Value = "COMP-THC-6117"
$ws.Range("J28").Value = "COMP-THC-6117"
$ws.Range("J29").Value = "THC-6118"
$ws.Range("J30").Value = "COMP-THC-6118"
$ws.Range("J31").Value = "THC-6118"
$ws.Range("J32").Value = "COMP-THC-6118"
$ws.Range("J33").Value = "COMP-THC-6119"
$ws.Range("J34").Value = "THC-6119"
$ws.Range("J35").Value = "COMP-THC-6119"
$ws.Range("J36").Value = "COMP-THC-6119"
$ws.Range("J37").Value = "THC-6119"
$ws.Range("J38").Value = "COMP-THC-6119"
$ws.Range("J39").Value = "COMP-THC-6119"
$ws.Range("J40").Value = "THC-6119"
$ws.Range("J41").Value = "THC-6120"
$ws.Range("J42").Value = "COMP-THC-6121"
$ws.Range("J43").Value = "COMP-THC-6121"
$ws.Range("J44").Value = "THC-6122"
$ws.Range("J45").Value = "THC-6122"
